# Apply the "cryptos" price-refresh update described by the commit's diff.
#
# The Price column (D) holds numbers formatted with "." as both the
# thousands separator and, in some rows, the decimal separator (e.g.
# "29.959.53", "1.001", "243.90") -- these are literal display strings,
# not real numeric values. Plain `.Value = "1.001"` assignment lets Excel
# reinterpret number-looking text as an actual number (dropping trailing
# zeros / losing the multi-dot grouping), so every such Price cell is
# written with a leading apostrophe to force a literal text value, exactly
# like typing `'1.001` into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.959.53'
$ws.Range("E2").Value = '  +0.51%  '

# Row 3
$ws.Range("D3").Value = '1.893.63'
$ws.Range("E3").Value = '  +0.11%  '

# Row 4
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").Value = '''0.7754'
$ws.Range("E5").Value = '  -0.32%  '

# Row 6
$ws.Range("D6").Value = '''243.90'
$ws.Range("E6").Value = '  +0.06%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("D8").Value = '''0.3135'
$ws.Range("E8").Value = '  +0.26%  '

# Row 9
$ws.Range("D9").Value = '''25.81'

# Row 10
$ws.Range("D10").Value = '''0.07269'
$ws.Range("E10").Value = '  +1.39%  '

# Row 11
$ws.Range("D11").Value = '''0.08685'
$ws.Range("E11").Value = '  +7.82%  '

# Row 12
$ws.Range("D12").Value = '2.007.00'
$ws.Range("E12").Value = '  +7.61%  '

# Row 13
$ws.Range("D13").Value = '''0.7744'
$ws.Range("E13").Value = '  +1.43%  '

# Row 14
$ws.Range("D14").Value = '''5.409'
$ws.Range("E14").Value = '  -0.83%  '

# Row 15
$ws.Range("D15").Value = '''94.47'
$ws.Range("E15").Value = '  +2.43%  '

# Row 16
$ws.Range("D16").Value = '''6.207'
$ws.Range("E16").Value = '  +0.81%  '

# Row 17
$ws.Range("D17").Value = '30.073.74'
$ws.Range("E17").Value = '  +0.97%  '

# Row 18
$ws.Range("D18").Value = '''13.93'
$ws.Range("E18").Value = '  -0.04%  '

# Row 19
$ws.Range("D19").Value = '''245.96'
$ws.Range("E19").Value = '  +1.06%  '

# Row 20
$ws.Range("D20").Value = '''0.000007882'
$ws.Range("E20").Value = '  +1.56%  '

# Row 21
$ws.Range("D21").Value = '2.267.77'
$ws.Range("E21").Value = '  +7.77%  '

# Row 22
$ws.Range("D22").Value = '''8.201'
$ws.Range("E22").Value = '  +1.26%  '

# Row 23
$ws.Range("D23").Value = '''1.001'
$ws.Range("E23").Value = '  +0.20%  '

# Row 24
$ws.Range("D24").Value = '''1.001'
$ws.Range("E24").Value = '  +0.06%  '

# Row 25
$ws.Range("D25").Value = '''0.1637'
$ws.Range("E25").Value = '  +1.33%  '

# Row 26
$ws.Range("D26").Value = '''9.510'
$ws.Range("E26").Value = '  +1.22%  '

# Row 27
$ws.Range("E27").Value = '  +1.06%  '

# Row 28
$ws.Range("D28").Value = '''18.83'
$ws.Range("E28").Value = '  +0.66%  '

# Row 29
$ws.Range("D29").Value = '''2.052'
$ws.Range("E29").Value = '  +0.23%  '

# Row 30
$ws.Range("D30").Value = '''1.431'
$ws.Range("E30").Value = '  +0.89%  '

# Row 31
$ws.Range("E31").Value = '  -0.24%  '

# Row 32
$ws.Range("D32").Value = '''4.519'
$ws.Range("E32").Value = '  +1.02%  '

# Row 33
$ws.Range("D33").Value = '''4.130'
$ws.Range("E33").Value = '  +0.69%  '

# Row 34
$ws.Range("D34").Value = '''0.05491'
$ws.Range("E34").Value = '  -0.79%  '

# Row 35
$ws.Range("D35").Value = '''1.250'
$ws.Range("E35").Value = '  -1.10%  '

# Row 36
$ws.Range("D36").Value = '''0.7554'
$ws.Range("E36").Value = '  +1.61%  '

# Row 37
$ws.Range("D37").Value = '''1.002'
$ws.Range("E37").Value = '  +0.57%  '

# Row 38
$ws.Range("D38").Value = '''2.688'
$ws.Range("E38").Value = '  +2.72%  '

# Row 39
$ws.Range("D39").Value = '''0.01962'
$ws.Range("E39").Value = '  +2.59%  '

# Row 40
$ws.Range("D40").Value = '''2.788'
$ws.Range("E40").Value = '  +0.06%  '

# Row 41
$ws.Range("D41").Value = '''0.4518'
$ws.Range("E41").Value = '  +2.24%  '

# Rows 42-43: Aave / Maker swap ranking order
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.108.72'
$ws.Range("E42").Value = '  -2.67%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '''73.70'
$ws.Range("E43").Value = '  -0.05%  '

# Row 44
$ws.Range("D44").Value = '''6.074'
$ws.Range("E44").Value = '  +3.76%  '

# Row 45
$ws.Range("D45").Value = '''0.8526'
$ws.Range("E45").Value = '  -0.08%  '

# Row 46
$ws.Range("D46").Value = '''0.9997'
$ws.Range("E46").Value = '  +0.00%  '

# Rows 47-49: Quant / RenderToken / RocketPoolETH rotate ranking order
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '2.188.76'
$ws.Range("E47").Value = '  +8.93%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '''103.20'
$ws.Range("E48").Value = '  -0.49%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''1.883'
$ws.Range("E49").Value = '  -0.05%  '

# Row 50
$ws.Range("D50").Value = '''7.614'
$ws.Range("E50").Value = '  +2.32%  '

# Row 51
$ws.Range("D51").Value = '''9.873'
$ws.Range("E51").Value = '  -0.25%  '
